# Automação para coletar temperatura
# Appends the latest sensor reading (Data/Hora, Temperatura, Umidade) as a new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "02/11/2024 20:55:19"
$ws.Range("B6").Value = "21º"

# Humidity is recorded as literal text (e.g. "88%"), not a numeric percentage,
# so force Text formatting before writing the value, then restore the default
# "Normal" style so no extra formatting is left behind on the cell.
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "88%"
$ws.Range("C6").Style = "Normal"
